# Added auto email and Hashing function
# - Adds a new "Monthly_STAT" worksheet summarizing visit counts.
# - Updates the roster on the "Sheet" tab (hashed id / name / visit count),
#   including a newly hashed entry for Steven Vargas.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the roster on the "Sheet" tab ---------------------------------
$ws.Range("A2").Value = "8A9AB340"
$ws.Range("B2").Value = "CBS"
$ws.Range("C2").Value = 69

$ws.Range("A3").Value = "CFD893A460"
$ws.Range("B3").Value = "Gary Tsai"
$ws.Range("C3").Value = 17
$ws.Rows.Item(3).RowHeight = 30

$ws.Range("A4").Value = "8FD8AAE4A0"
$ws.Range("B4").Value = "Steven Vargas"
$ws.Range("C4").Value = 2

$ws.Range("H17").Select()

# --- Add the Monthly_STAT worksheet ---------------------------------------
$stat = $wb.Worksheets.Add($null, $ws)
$stat.Name = "Monthly_STAT"

$stat.Range("A3").Value = "Total student visited this month: "
$stat.Range("C3").Formula = "=SUM(sheet!C2:C4)"
$stat.Rows.Item(3).RowHeight = 23

$labelRange = $stat.Range("A3:C3")
$labelRange.Font.Italic = $true
$labelRange.Font.Size = 18
$labelRange.Font.Color = 0

$stat.Range("C3").Select()
